$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slang")

# Add new block of data mirroring the existing ones (rows 1-2 and 5-6)
# Row 8: label + bold "User"/"Password" headers
$ws.Range("A8").Value = "mobileTest"

$ws.Range("B8").Value = "User"
$ws.Range("B8").Font.Bold = $true

$ws.Range("C8").Value = "Password"
$ws.Range("C8").Font.Bold = $true

# Row 9: hyperlink-styled sample values
$ws.Range("B9").Value = "sampleusername"
$ws.Range("B9").Style = "Hyperlink"

$ws.Range("C9").Value = "samplePwd"
$ws.Range("C9").Style = "Hyperlink"

# Move the active selection, matching the saved view state
$ws.Range("A14").Select()
